# 1) Mark the run that holds the inline drawing as "do not spell/grammar
#    check" -> adds <w:rPr><w:noProof/></w:rPr> to that run, matching the
#    diff's first hunk.
$d = $word.ActiveDocument
$imgPara = $d.Paragraphs.Item(1)
$imgPara.Range.NoProofing = $true

# 2) Append the new paragraphs (blank line, red/bold/underlined heading,
#    and the three SQL command paragraphs) after the existing last
#    paragraph, matching the diff's second hunk.
$end = $d.Content
$end.Collapse(0)

$frag = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'/>" `
    + "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" `
    + "<w:pPr><w:rPr><w:b/><w:bCs/><w:color w:val='FF0000'/><w:u w:val='single'/></w:rPr></w:pPr>" `
    + "<w:r><w:rPr><w:b/><w:bCs/><w:color w:val='FF0000'/><w:u w:val='single'/></w:rPr>" `
    + "<w:t xml:space='preserve'>Este es el código para </w:t></w:r>" `
    + "<w:proofErr w:type='gramStart'/>" `
    + "<w:r><w:rPr><w:b/><w:bCs/><w:color w:val='FF0000'/><w:u w:val='single'/></w:rPr>" `
    + "<w:t>hacerlo :</w:t></w:r>" `
    + "<w:proofErr w:type='gramEnd'/>" `
    + "</w:p>" `
    + "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" `
    + "<w:r><w:t>CREATE USER 'nascor04'@'localhost' IDENTIFIED BY 'Nasc0r2020!';</w:t></w:r>" `
    + "</w:p>" `
    + "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" `
    + "<w:r><w:t>GRANT ALL PRIVILEGES ON nascor04_</w:t></w:r>" `
    + "<w:proofErr w:type='gramStart'/>" `
    + "<w:r><w:t>bddCurso.*</w:t></w:r>" `
    + "<w:proofErr w:type='gramEnd'/>" `
    + "<w:r><w:t xml:space='preserve'> TO 'nascor04'@'localhost';</w:t></w:r>" `
    + "</w:p>" `
    + "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" `
    + "<w:r><w:t>FLUSH PRIVILEGES;</w:t></w:r>" `
    + "</w:p>"

[void]$end.InsertXML($frag)

Write-Output "Edit applied. Paragraph count: $($d.Paragraphs.Count)"
